$wb = $excel.ActiveWorkbook

# --- Sheet "string" ---
$ws1 = $wb.Worksheets.Item("string")
$ws1.Range("A1").Value = "turn"
$ws1.Range("B1").Value = "varname"

# --- Sheet "numeric" ---
$ws2 = $wb.Worksheets.Item("numeric")
$ws2.Range("A1").Value = "turn"
$ws2.Range("B1").Value = "varname"

# --- Sheet "drop" ---
$ws3 = $wb.Worksheets.Item("drop")
$ws3.Range("A1").Value = "turn"
$ws3.Range("B1").Value = "n_obs"
$ws3.Range("C1").Value = "initials"
$ws3.Range("D1").Value = "notes"

# Extend formatting of the data row (row 2) to the new column D,
# matching the style already used by A2:C2
$ws3.Range("C2").Copy() | Out-Null
$ws3.Range("D2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
